# Final updates: IP block fix, timestamped report, QR rules
# Add two new student records (roll no. + name) below the existing list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows
$ws.Range("A5").Value = "23BCA004"
$ws.Range("B5").Value = "Drishti paras"
$ws.Range("A6").Value = "23BCA005"
$ws.Range("B6").Value = "Ishan sharma"

# Match the formatting (style + row height) of the preceding data row
$ws.Range("A5:B6").HorizontalAlignment = 1
$ws.Range("5:6").RowHeight = $ws.Range("4:4").RowHeight
